$d = $word.ActiveDocument

$old = "2018年キャンペーン期間 (対象：ペルセウス)：、10月30日〜11月8日、11月29日〜12月8日"
$new = "年キャンペーン期間 対象：Leo: 4月14〜23日、5月14〜23日"

for ($i = 0; $i -lt 10; $i++) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        break
    }
    $rng.Delete()
    $rng.InsertAfter($new)
}
